$wb = $excel.ActiveWorkbook
Write-Host $wb.Windows.Count
$win = $excel.ActiveWindow
Write-Host $win.WindowState
